$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-74)
# from serial date 45190 (2023-09-21) to 45192 (2023-09-23)
for ($r = 2; $r -le 74; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -eq 45190) {
        $cell.Value = 45192
    }
}
